$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1586
$ws.Range("E2").Value = 52
$ws.Range("F2").Value = 52
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = -1
$ws.Range("I2").Value = -1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1602
$ws.Range("L2").Value = 1158
$ws.Range("M2").Value = 445
$ws.Range("N2").Value = 437
$ws.Range("O2").Value = 8
$ws.Range("P2").Value = 197
$ws.Range("Q2").Value = 342
$ws.Range("R2").Value = -36
$ws.Range("S2").Value = -134
$ws.Range("T2").Value = 36
$ws.Range("U2").Value = 305
$ws.Range("V2").Value = 750
$ws.Range("W2").Value = 3.26
$ws.Range("X2").Value = -0.07000000000000001
$ws.Range("Y2").Value = -0.25
$ws.Range("Z2").Value = -0.07000000000000001
$ws.Range("AA2").Value = 260.18
$ws.Range("AB2").Value = 58.91
$ws.Range("AC2").Value = -3
$ws.Range("AD2").Value = -482.78
$ws.Range("AE2").Value = 1300
$ws.Range("AF2").Value = 1.04
$ws.Range("AG2").Value = 60
$ws.Range("AH2").Value = 4.44
$ws.Range("AI2").Value = -1832.49
$ws.Range("AJ2").Value = 39353308

# Row 3
$ws.Range("D3").Value = 1590
$ws.Range("E3").Value = 69
$ws.Range("F3").Value = 69
$ws.Range("G3").Value = 213
$ws.Range("H3").Value = 145
$ws.Range("I3").Value = 145
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1654
$ws.Range("L3").Value = 1075
$ws.Range("M3").Value = 580
$ws.Range("N3").Value = 571
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 197
$ws.Range("Q3").Value = 73
$ws.Range("R3").Value = -142
$ws.Range("S3").Value = 18
$ws.Range("T3").Value = 77
$ws.Range("U3").Value = -4
$ws.Range("V3").Value = 799
$ws.Range("W3").Value = 4.31
$ws.Range("X3").Value = 9.1
$ws.Range("Y3").Value = 28.69
$ws.Range("Z3").Value = 8.890000000000001
$ws.Range("AA3").Value = 185.41
$ws.Range("AB3").Value = 121.64
$ws.Range("AC3").Value = 367
$ws.Range("AD3").Value = 5.09
$ws.Range("AE3").Value = 1700
$ws.Range("AF3").Value = 1.1
$ws.Range("AG3").Value = 60
$ws.Range("AH3").Value = 3.21
$ws.Range("AI3").Value = 13.94
$ws.Range("AJ3").Value = 39353308

# Row 4
$ws.Range("D4").Value = 1687
$ws.Range("E4").Value = 121
$ws.Range("F4").Value = 121
$ws.Range("G4").Value = -262
$ws.Range("H4").Value = -205
$ws.Range("I4").Value = -204
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 1789
$ws.Range("L4").Value = 1444
$ws.Range("M4").Value = 345
$ws.Range("N4").Value = 338
$ws.Range("O4").Value = 7
$ws.Range("P4").Value = 197
$ws.Range("Q4").Value = 80
$ws.Range("R4").Value = -254
$ws.Range("S4").Value = 116
$ws.Range("T4").Value = 178
$ws.Range("U4").Value = -98
$ws.Range("V4").Value = 925
$ws.Range("W4").Value = 7.19
$ws.Range("X4").Value = -12.17
$ws.Range("Y4").Value = -44.99
$ws.Range("Z4").Value = -11.92
$ws.Range("AA4").Value = 418.85
$ws.Range("AB4").Value = 14.12
$ws.Range("AC4").Value = -519
$ws.Range("AD4").Value = -5.75
$ws.Range("AE4").Value = 1004
$ws.Range("AF4").Value = 2.97
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 39353308

# Row 5
$ws.Range("D5").Value = 1804
$ws.Range("E5").Value = 131
$ws.Range("F5").Value = 131
$ws.Range("G5").Value = 325
$ws.Range("H5").Value = 186
$ws.Range("I5").Value = 187
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1961
$ws.Range("L5").Value = 1385
$ws.Range("M5").Value = 576
$ws.Range("N5").Value = 569
$ws.Range("O5").Value = 7
$ws.Range("P5").Value = 214
$ws.Range("Q5").Value = -221
$ws.Range("R5").Value = 361
$ws.Range("S5").Value = -135
$ws.Range("T5").Value = 43
$ws.Range("U5").Value = -264
$ws.Range("V5").Value = 706
$ws.Range("W5").Value = 7.29
$ws.Range("X5").Value = 10.34
$ws.Range("Y5").Value = 41.19
$ws.Range("Z5").Value = 9.949999999999999
$ws.Range("AA5").Value = 240.45
$ws.Range("AB5").Value = 133.08
$ws.Range("AC5").Value = 446
$ws.Range("AD5").Value = 5.56
$ws.Range("AE5").Value = 1540
$ws.Range("AF5").Value = 1.61
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 42710111

# Row 6
$ws.Range("D6").Value = 2054
$ws.Range("E6").Value = 135
$ws.Range("F6").Value = 135
$ws.Range("G6").Value = 272
$ws.Range("H6").Value = 203
$ws.Range("I6").Value = 202
$ws.Range("K6").Value = 2048
$ws.Range("L6").Value = 1253
$ws.Range("M6").Value = 795
$ws.Range("N6").Value = 788
$ws.Range("P6").Value = 216
$ws.Range("Q6").Value = 174
$ws.Range("R6").Value = -120
$ws.Range("S6").Value = 30
$ws.Range("T6").Value = 32
$ws.Range("U6").Value = 142
$ws.Range("V6").Value = 722
$ws.Range("W6").Value = 6.58
$ws.Range("X6").Value = 9.859999999999999
$ws.Range("Y6").Value = 29.84
$ws.Range("Z6").Value = 10.11
$ws.Range("AA6").Value = 157.72
$ws.Range("AB6").Value = 230.16
$ws.Range("AC6").Value = 469
$ws.Range("AD6").Value = 4.98
$ws.Range("AE6").Value = 2062
$ws.Range("AF6").Value = 1.13
$ws.Range("AG6").Value = 80
$ws.Range("AH6").Value = 3.43
$ws.Range("AI6").Value = 15.1
$ws.Range("AJ6").Value = 43226540

# Clear rows 7-9 data cells (D:AI), keep A/B/C
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
